$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (tab) name from "New Data" to "New data"
$ws.Name = "New data"

# Clear the old A1:C6 range contents, since the new table shape is A1:F5
$ws.Range("A1:C6").ClearContents()

# Header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Branch"
$ws.Range("C1").Value = "Sem1"
$ws.Range("D1").Value = "Sem2"
$ws.Range("E1").Value = "title1"
$ws.Range("F1").Value = "title"

# Row 2 - pradip
$ws.Range("A2").Value = "pradip"
$ws.Range("B2").Value = "EC"
$ws.Range("C2").Value = 9.2
$ws.Range("D2").Value = 6.9
$ws.Range("E2").Value = "manager"
$ws.Range("F2").Value = "president"

# Row 3 - darshan
$ws.Range("A3").Value = "darshan"
$ws.Range("B3").Value = "IT"
$ws.Range("C3").Value = 8.8
$ws.Range("D3").Value = 7.9
$ws.Range("E3").Value = "manager"
$ws.Range("F3").Value = "president"

# Row 4 - jay
$ws.Range("A4").Value = "jay"
$ws.Range("B4").Value = "IC"
$ws.Range("C4").Value = 6.9
$ws.Range("D4").Value = 9.2
$ws.Range("E4").Value = "manager"
$ws.Range("F4").Value = "president"

# Row 5 - vishal
$ws.Range("A5").Value = "vishal"
$ws.Range("B5").Value = "EC"
$ws.Range("C5").Value = 8.4
$ws.Range("D5").Value = 8.8
$ws.Range("E5").Value = "manager"
$ws.Range("F5").Value = "president"
